# Apply the changes described by the diff:
# 1. Clear the stray empty cell B2 on the "ODI Batting" sheet.
# 2. Add a new worksheet "ODI Batting Extra" at the end of the workbook with
#    headers + three rows of data.

$wb = $excel.ActiveWorkbook

# --- 1. Clear B2 on "ODI Batting" -------------------------------------------------
$wsBatting = $wb.Worksheets.Item("ODI Batting")
$wsBatting.Range("B2").ClearContents()

# --- 2. Add the new "ODI Batting Extra" sheet after the last existing sheet -------
$lastIndex = $wb.Worksheets.Count
$wsExtra = $wb.Worksheets.Add([Type]::Missing, $wb.Worksheets.Item($lastIndex))
$wsExtra.Name = "ODI Batting Extra"

# Header row (bold / bordered / centered style, matching the other sheets' header look)
$headers = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($col = 1; $col -le $headers.Length; $col++) {
    $wsExtra.Cells.Item(1, $col).Value = $headers[$col - 1]
}
# Copy the header formatting (bold font, thin border, centered) from the
# "ODI Batting" sheet's own header row so the new sheet matches the workbook's look.
$wsBatting.Range("A1").Copy()
$wsExtra.Range("A1:F1").PasteSpecial(-4122)  # xlPasteFormats

# Data rows: MATCH_CODE, BATTING_POSITION, NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH
$rows = @(
    @("4472", 6, "", "", "", "NO"),
    @("4473", 6, "2", "0", "6.88%", "NO"),
    @("4476", 6, "0", "0", "0.90%", "NO")
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowData = $rows[$r]
    $excelRow = $r + 2
    $wsExtra.Cells.Item($excelRow, 1).Value = "'" + $rowData[0]
    $wsExtra.Cells.Item($excelRow, 2).Value = $rowData[1]
    $wsExtra.Cells.Item($excelRow, 3).Value = "'" + $rowData[2]
    $wsExtra.Cells.Item($excelRow, 4).Value = "'" + $rowData[3]
    $wsExtra.Cells.Item($excelRow, 5).Value = "'" + $rowData[4]
    $wsExtra.Cells.Item($excelRow, 6).Value = $rowData[5]
}
